# Update values in Sheet1 ("Name of Algo" update per commit message).
# The source data changed slightly, producing new imputed/measured values
# in columns A and C for a number of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = -22.22429999999999
$ws.Range("A14").Value = -21.72920000000001
$ws.Range("A21").Value = -20.09199999999998
$ws.Range("C22").Value = -11.6696
$ws.Range("A23").Value = -20.32659999999997
$ws.Range("C24").Value = -13.46949999999999
$ws.Range("A25").Value = -21.59359999999998
$ws.Range("A26").Value = -21.06689999999996
$ws.Range("C28").Value = -13.72509999999999
$ws.Range("A29").Value = -20.83009999999997
$ws.Range("C36").Value = -12.0609
$ws.Range("C45").Value = -13.64699999999999
$ws.Range("C48").Value = -12.3703
$ws.Range("C49").Value = -13.63919999999998
$ws.Range("C52").Value = -10.8126
$ws.Range("A53").Value = -21.61040000000001
$ws.Range("C53").Value = -10.93300000000001
$ws.Range("C54").Value = -13.4284
$ws.Range("A57").Value = -22.3259
$ws.Range("A59").Value = -21.9489
$ws.Range("A69").Value = -21.5699
$ws.Range("C70").Value = -12.6393
$ws.Range("A79").Value = -20.2809
$ws.Range("A83").Value = -22.03830000000001
$ws.Range("C86").Value = -13.86499999999999
$ws.Range("C87").Value = -12.7407
$ws.Range("C89").Value = -13.1436
$ws.Range("A91").Value = -20.38969999999998
$ws.Range("A93").Value = -21.20060000000001
$ws.Range("C101").Value = -13.3159
$ws.Range("A103").Value = -21.78929999999999
